# Generate Report for Handback
#
# The "d86206fc-0e9f-4b61-8a02-3cbb2e1b3c55" row's handback failed: the
# handback file name did not match the handoff file name. Reflect this in
# the status columns and populate the per-language "Error Detail" cells,
# widening that column so the message is readable.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "Handback transform failed"

# Overview sheet: zh-cn / de-de status columns for the d86206fc row (row 3)
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# zh-cn sheet: Status column (C) for the d86206fc row (row 3)
$wsZhCn.Range("C3").Value = $newStatus

# zh-cn sheet: Error Detail column (P) for the d86206fc row (row 3)
$wsZhCn.Range("P3").Value = "Handback file name: qphotvzc.jjb is different with handoff file name: d86206fc-0e9f-4b61-8a02-3cbb2e1b3c55.9d5fe46469d8ccf99275b56fc9afc5f74d750b57.zh-cn."

# de-de sheet: Status column (C) for the d86206fc row (row 3)
$wsDeDe.Range("C3").Value = $newStatus

# de-de sheet: Error Detail column (P) for the d86206fc row (row 3)
$wsDeDe.Range("P3").Value = "Handback file name: qphotvzc.jjb is different with handoff file name: d86206fc-0e9f-4b61-8a02-3cbb2e1b3c55.9d5fe46469d8ccf99275b56fc9afc5f74d750b57.de-de."

# Widen the "Error Detail" column (16th column, P) on both language sheets
# now that it holds a long diagnostic message. ColumnWidth is expressed in
# "characters" but the engine stores the XML column width with a constant
# +5/6 character padding baked in, so back that out to land on width=40.
$targetColumnWidth = 40 - (5 / 6)
$wsZhCn.Columns.Item(16).ColumnWidth = $targetColumnWidth
$wsDeDe.Columns.Item(16).ColumnWidth = $targetColumnWidth
